$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "2023-10-22 21:36:58 7 question_7_8530875"
$ws.Range("B39").Value = "7 question"

# Column C holds the numeric-looking "7" as text (matching the rest of the
# sheet, e.g. rows that store "1", "2", "3" ... as shared strings rather
# than numbers). Briefly flip the cell to a text format so the value is
# committed as a string, then restore the default "General" format so the
# cell's style stays identical to its neighbours (style index 0).
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "7"
$ws.Range("C39").NumberFormat = "General"
